$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.5587383333333333
$ws.Range("H2").Value = 1.676215
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.676435666666666
$ws.Range("N2").Value = 11.029307
$ws.Range("O2").Value = 0.05596928005870617
$ws.Range("P2").Value = 0.07888124434163156
$ws.Range("Q2").Value = 2.054165537000555
$ws.Range("R2").Value = 18.487489833005
$ws.Range("S2").Value = 0.05596928005870617
$ws.Range("T2").Value = 0.07888124434163156

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.5587383333333333
$ws.Range("H3").Value = 1.676215
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.794922333333334
$ws.Range("N3").Value = 11.384767
$ws.Range("O3").Value = 0.05777309604548284
$ws.Range("P3").Value = 0.08142348268114613
$ws.Range("Q3").Value = 2.120368579656111
$ws.Range("R3").Value = 19.083317216905
$ws.Range("S3").Value = 0.05777309604548284
$ws.Range("T3").Value = 0.08142348268114613

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.5587383333333333
$ws.Range("H4").Value = 1.676215
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.477427
$ws.Range("N4").Value = 1.432281
$ws.Range("O4").Value = 0.007268247806663079
$ws.Range("P4").Value = 0.01024362704990227
$ws.Range("Q4").Value = 0.2667567662683333
$ws.Range("R4").Value = 2.400810896415
$ws.Range("S4").Value = 0.007268247806663079
$ws.Range("T4").Value = 0.01024362704990227

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.5587383333333333
$ws.Range("H5").Value = 1.676215
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.4995396666666667
$ws.Range("N5").Value = 1.498619
$ws.Range("O5").Value = 0.007604886373395734
$ws.Range("P5").Value = 0.01071807426468513
$ws.Range("Q5").Value = 0.2791119607872222
$ws.Range("R5").Value = 2.512007647085
$ws.Range("S5").Value = 0.007604886373395734
$ws.Range("T5").Value = 0.01071807426468513

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.5587383333333333
$ws.Range("H6").Value = 1.676215
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 57.23834600000001
$ws.Range("N6").Value = 114.476692
$ws.Range("O6").Value = 0.8713844897157522
$ws.Range("P6").Value = 0.8187335716626348
$ws.Range("Q6").Value = 31.98125804679667
$ws.Range("R6").Value = 191.88754828078
$ws.Range("S6").Value = 0.8713844897157522
$ws.Range("T6").Value = 0.8187335716626348

